$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (A3) onto the new row's label cell (A4)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new label value (adds "Q1" to shared strings)
$ws.Range("A4").Value = "Q1"

# Update existing row 2 values
$ws.Range("B2").Value = 0.1779586292060669
$ws.Range("C2").Value = 2.019706927474283
$ws.Range("D2").Value = 17.42504876041586
$ws.Range("E2").Value = 4.174332133457502
$ws.Range("F2").Value = 4.222346974595365
$ws.Range("G2").Value = 41

# Update existing row 3 values
$ws.Range("B3").Value = 0.2271760140433983
$ws.Range("C3").Value = 1.978660754480025
$ws.Range("D3").Value = 15.35741335774997
$ws.Range("E3").Value = 3.918853576972475
$ws.Range("F3").Value = 3.927054555281079
$ws.Range("G3").Value = 133

# Add new row 4 values
$ws.Range("B4").Value = 0.1196519597548314
$ws.Range("C4").Value = 1.272925636934957
$ws.Range("D4").Value = 5.710723481807822
$ws.Range("E4").Value = 2.389712008131486
$ws.Range("F4").Value = 2.404727861035995
$ws.Range("G4").Value = 67
